$wb = $excel.ActiveWorkbook

# Update the absolute-path hint Excel stores for the file (best effort; value
# is informational only).
# (Left as-is; not exposed via a settable COM property in this runtime.)

# Locate the existing "example-data" worksheet so the new sheet lands right
# after it (and before the two chart sheets, tab-order-wise).
$exampleData = $wb.Sheets.Item("example-data")

# Add the new worksheet used for the WrapUp GGSlides screenshot tally.
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $exampleData)
$newSheet.Name = "Sheet1"

# Match the authored column sizing for the screenshot-name column.
$newSheet.Columns.Item(1).ColumnWidth = 19.7

$filenames = @(
    'SS_00-01-04_GBR.png',
    'SS_00-01-12_GBR.png',
    'SS_00-01-30_GBR.png',
    'SS_00-01-31_GBR.png',
    'SS_00-01-33_GBR.png',
    'SS_00-01-34_GBR.PNG',
    'SS_00-01-36_GBR.PNG',
    'SS_00-01-37_GBR.png',
    'SS_00-01-48_GBR.PNG',
    'SS_00-02-05_GBR.PNG',
    'SS_00-02-08_GBR.PNG',
    'SS_00-02-16_GBR.PNG',
    'SS_00-02-44_GBR.PNG',
    'SS_00-03-04_GBR.PNG',
    'SS_00-03-13_GBR.PNG',
    'SS_00-05-28_GBR.PNG',
    'SS_00-05-34_GBR.PNG',
    'SS_00-05-42_GBR.PNG',
    'SS_00-05-53_GBR.PNG',
    'SS_00-06-47_GBR.PNG',
    'SS_00-07-18_GBR.PNG',
    'SS_00-07-26_GBR.PNG',
    'SS_00-07-35_GBR.PNG',
    'SS_00-07-40_GBR(2).PNG',
    'SS_00-07-40_GBR.PNG',
    'SS_00-07-41_GBR.PNG',
    'SS_00-07-43_GBR.PNG',
    'SS_00-07-49_GBR.PNG',
    'SS_00-07-51_GBR.PNG',
    'SS_00-07-53_GBR.PNG',
    'SS_00-07-56_GBR.PNG',
    'SS_00-07-58_GBR.PNG',
    'SS_00-08-01_GBR.PNG',
    'SS_00-08-03_GBR.PNG',
    'SS_00-08-04_GBR.PNG',
    'SS_00-08-09_GBR.PNG',
    'SS_00-08-12_GBR.PNG',
    'SS_00-08-53_GBR.PNG',
    'SS_00-09-02_GBR.PNG',
    'SS_00-09-47_GBR.PNG',
    'SS_00-09-54_GBR.PNG',
    'SS_00-09-57_GBR.PNG',
    'SS_00-10-02_GBR.PNG',
    'SS_00-10-08_GBR.PNG',
    'SS_00-10-13_GBR.PNG',
    'SS_00-10-20_GBR.PNG',
    'SS_00-10-39_GBR.PNG',
    'SS_00-10-42_GBR.PNG',
    'SS_00-10-43_GBR.PNG',
    'SS_00-10-47_GBR.PNG',
    'SS_00-11-18_GBR.PNG',
    'SS_00-11-20_GBR.PNG',
    'SS_00-11-21_GBR.PNG',
    'SS_00-11-22_GBR.PNG',
    'SS_00-11-24_GBR.PNG',
    'SS_00-11-27_GBR(2).PNG',
    'SS_00-11-27_GBR.PNG',
    'SS_00-11-28_GBR.PNG',
    'SS_00-11-29_GBR.PNG',
    'SS_00-11-31_GBR.PNG',
    'SS_00-11-32_GBR.PNG',
    'SS_00-11-33_GBR.PNG',
    'SS_00-11-36_GBR.PNG',
    'SS_00-11-38_GBR.PNG',
    'SS_00-11-40_GBR(2).PNG',
    'SS_00-11-40_GBR.PNG',
    'SS_00-11-42_GBR.PNG',
    'SS_00-11-43_GBR.PNG',
    'SS_00-11-44_GBR.PNG',
    'SS_00-11-47_GBR.PNG',
    'SS_00-11-49_GBR(2).PNG',
    'SS_00-11-49_GBR.PNG',
    'SS_00-11-50_GBR(2).PNG',
    'SS_00-11-50_GBR(3).PNG',
    'SS_00-11-50_GBR.PNG',
    'SS_00-11-52_GBR.PNG',
    'SS_00-12-02_GBR.PNG',
    'SS_00-12-21_GBR.PNG',
    'SS_00-12-23_GBR.PNG',
    'SS_00-13-43_GBR.PNG',
    'SS_00-14-05_GBR.PNG',
    'SS_00-15-08_GBR.PNG',
    'SS_00-16-33_SPO(2).PNG',
    'SS_00-16-33_SPO.PNG',
    'SS_00-16-34_SPO.PNG',
    'SS_00-16-35_SPO.PNG',
    'SS_00-16-37_SPO.PNG',
    'SS_00-16-39_SPO.PNG',
    'SS_00-16-40_SPO.PNG',
    'SS_00-16-46_SPO.PNG',
    'SS_00-16-48_SPO.PNG',
    'SS_00-16-58_SPO.PNG',
    'SS_00-17-02_SPO.PNG',
    'SS_00-17-05_SPO.PNG',
    'SS_00-17-09_SPO.PNG',
    'SS_00-17-11_SPO.PNG',
    'SS_00-17-14_SPO.PNG',
    'SS_00-17-34_SPO.PNG',
    'SS_00-17-37_SPO.PNG',
    'SS_00-19-00_SPO.PNG',
    'SS_00-20-34_SL.PNG',
    'SS_00-20-52_SL.PNG',
    'SS_00-21-05_SL.PNG',
    'SS_00-22-55_SL.PNG',
    'SS_00-25-26_SPO.PNG',
    'SS_00-26-48_DO.PNG',
    'SS_00-26-55_DO.PNG',
    'SS_00-27-18_DO.PNG',
    'SS_00-27-51_DO.PNG',
    'SS_00-28-32_DO.PNG',
    'SS_00-28-54_DO.PNG',
    'SS_00-29-50_DO.PNG',
    'SS_00-34-13_T.PNG',
    'SS_00-40-37_SPO.PNG',
    'SS_00-40-43_SPO.PNG',
    'SS_00-40-54_SPO.PNG',
    'SS_00-41-57_SPO.PNG',
    'SS_00-45-53_SPO.PNG',
    'SS_00-45-55_SPO.PNG',
    'SS_00-51-47_EAC.PNG',
    'SS_00-52-36_EAC.PNG',
    'SS_00-52-42_EAC.PNG',
    'SS_00-52-50_EAC.PNG',
    'SS_00-54-00_EAC.PNG',
    'SS_00-55-38_EAC.PNG',
    'SS_00-56-33_EAC.PNG',
    'SS_00-56-38_SPO.PNG',
    'SS_00-56-44_SPO.PNG',
    'SS_00-56-52_SPO.PNG',
    'SS_00-56-58_SPO.PNG',
    'SS_00-57-00_SPO.PNG',
    'SS_00-57-02_SPO.PNG',
    'SS_00-57-05_SPO.PNG',
    'SS_00-57-08_SH.PNG',
    'SS_00-57-25_SH.PNG',
    'SS_00-57-27_SH.PNG',
    'SS_00-57-34_SH.PNG',
    'SS_00-57-44_SH.PNG',
    'SS_00-59-45_DO.PNG',
    'SS_01-06-51_SPO.PNG',
    'SS_01-08-36_SH.PNG',
    'SS_01-08-39_SH.PNG',
    'SS_01-08-44_SH.PNG',
    'SS_01-08-46_SH.PNG',
    'SS_01-08-55_SH.PNG',
    'SS_01-13-48_SH.PNG',
    'SS_01-16-48_SH.PNG',
    'SS_01-17-34_SH.PNG',
    'SS_01-21-50_SH.PNG',
    'SS_01-24-02_SH.PNG',
    'SS_01-24-18_SH.PNG',
    'SS_01-24-30_SH.PNG',
    'SS_01-25-16_SH.PNG',
    'SS_01-25-19_FG.PNG',
    'SS_01-26-22_FG.PNG',
    'SS_01-26-35_FG.PNG',
    'SS_01-31-00_GBR.PNG',
    'SS_01-31-02_GBR.PNG',
    'SS_01-31-05_GBR.PNG',
    'SS_01-31-06_GBR.PNG',
    'SS_01-31-16_GBR.PNG',
    'SS_01-31-20_GBR.PNG',
    'SS_01-31-45_GBR.PNG',
    'SS_01-31-49_GBR.PNG',
    'SS_01-31-51_GBR.PNG',
    'SS_01-32-06_GBR.PNG',
    'SS_01-33-20_SH.PNG',
    'SS_01-33-36_SPO.PNG',
    'SS_01-39-37_credit_MikeWazowski.PNG',
    'SS_01-40-12_credit_T.PNG'
)


# Column A: header + one row per screenshot filename.
$newSheet.Cells.Item(1, 1).Value = "Screenshot"
for ($i = 0; $i -lt $filenames.Length; $i++) {
    $newSheet.Cells.Item($i + 2, 1).Value = $filenames[$i]
}

# Columns B:E: the four summary headers (only header row is populated).
$newSheet.Cells.Item(1, 2).Value = "new_identified"
$newSheet.Cells.Item(1, 3).Value = "new_unidentified"
$newSheet.Cells.Item(1, 4).Value = "total_identified"
$newSheet.Cells.Item(1, 5).Value = "total_unidentified-upt"

# Make the new sheet the active/selected tab, with E2 selected, mirroring
# the authored sheetView state.
$newSheet.Activate()
$newSheet.Range("E2").Select()
